$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 2256.3845
$ws.Range("I131").Value = 404.7143
$ws.Range("J131").Value = 4416.6665
$ws.Range("K131").Value = 1214.1429
$ws.Range("L131").Value = 13249.9995
$ws.Range("M131").Value = 3825.8571
$ws.Range("N131").Value = -23329.9995
$ws.Range("H137").Value = 2106.5881
$ws.Range("I137").Value = 1397.875
$ws.Range("J137").Value = 3807.5
$ws.Range("K137").Value = 4193.625
$ws.Range("L137").Value = 11422.5
$ws.Range("M137").Value = -1643.625
$ws.Range("N137").Value = -16522.5
$ws.Range("H141").Value = 7734.161
$ws.Range("I141").Value = 8954.16
$ws.Range("J141").Value = 2650.8333
$ws.Range("K141").Value = 26862.48
$ws.Range("L141").Value = 7952.499899999999
$ws.Range("M141").Value = -21682.48
$ws.Range("N141").Value = -18312.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4853.48
$ws.Range("I74").Value = 5563.1875
$ws.Range("J74").Value = 3591.7778
$ws.Range("K74").Value = 5563.1875
$ws.Range("L74").Value = 3591.7778
$ws.Range("M74").Value = -4689.1875
$ws.Range("N74").Value = -5339.7778
$ws.Range("H77").Value = 4853.48
$ws.Range("I77").Value = 5563.1875
$ws.Range("J77").Value = 3591.7778
$ws.Range("K77").Value = 27815.9375
$ws.Range("L77").Value = 17958.889
$ws.Range("M77").Value = -23447.9375
$ws.Range("N77").Value = -26694.889
$ws.Range("H88").Value = 8337170.5
$ws.Range("I88").Value = 9527481
$ws.Range("K88").Value = 9527481
$ws.Range("M88").Value = -9527075
$ws.Range("H91").Value = 8337170.5
$ws.Range("I91").Value = 9527481
$ws.Range("K91").Value = 9527481
$ws.Range("M91").Value = -9526077

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("N9").ClearContents()
$ws.Range("H86").Value = 2092.3635
$ws.Range("I86").Value = 1828.3334
$ws.Range("J86").Value = 2658.1428
$ws.Range("K86").Value = 1828.3334
$ws.Range("L86").Value = 2658.1428
$ws.Range("M86").Value = -705.3334
$ws.Range("N86").Value = -4904.1428
$ws.Range("H89").Value = 2092.3635
$ws.Range("I89").Value = 1828.3334
$ws.Range("J89").Value = 2658.1428
$ws.Range("K89").Value = 9141.666999999999
$ws.Range("L89").Value = 13290.714
$ws.Range("M89").Value = -3525.666999999999
$ws.Range("N89").Value = -24522.714
$ws.Range("H94").Value = 2250
$ws.Range("J94").Value = 1500
$ws.Range("L94").Value = 1500
$ws.Range("N94").Value = -2402
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H134").Value = 2541.8125
$ws.Range("I134").Value = 1651.091
$ws.Range("J134").Value = 4501.4
$ws.Range("K134").Value = 4953.272999999999
$ws.Range("L134").Value = 13504.2
$ws.Range("M134").Value = -2418.272999999999
$ws.Range("N134").Value = -18574.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6174044.5
$ws.Range("J16").Value = 1355.2222
$ws.Range("L16").Value = 1355.2222
$ws.Range("N16").Value = -1929.2222
$ws.Range("H31").Value = 12197145
$ws.Range("I31").Value = 801.48
$ws.Range("J31").Value = 31253930
$ws.Range("K31").Value = 801.48
$ws.Range("L31").Value = 31253930
$ws.Range("M31").Value = -506.48
$ws.Range("N31").Value = -31254520
$ws.Range("H34").Value = 12197145
$ws.Range("I34").Value = 801.48
$ws.Range("J34").Value = 31253930
$ws.Range("K34").Value = 801.48
$ws.Range("L34").Value = 31253930
$ws.Range("M34").Value = -599.48
$ws.Range("N34").Value = -31254334
$ws.Range("H58").Value = 1470.1333
$ws.Range("I58").Value = 1426.9595
$ws.Range("K58").Value = 1426.9595
$ws.Range("M58").Value = -1223.9595
$ws.Range("H62").Value = 7000
$ws.Range("J62").Value = 7000
$ws.Range("L62").Value = 7000
$ws.Range("N62").Value = -8248
$ws.Range("H65").Value = 7000
$ws.Range("J65").Value = 7000
$ws.Range("L65").Value = 35000
$ws.Range("N65").Value = -41240
$ws.Range("H105").Value = 2392.8
$ws.Range("I105").Value = 2003.6
$ws.Range("J105").Value = 2782
$ws.Range("K105").Value = 2003.6
$ws.Range("L105").Value = 2782
$ws.Range("M105").Value = -256.5999999999999
$ws.Range("N105").Value = -6276
$ws.Range("H113").Value = 6174044.5
$ws.Range("J113").Value = 1355.2222
$ws.Range("L113").Value = 1355.2222
$ws.Range("N113").Value = -5695.2222
$ws.Range("H134").Value = 4221.472
$ws.Range("I134").Value = 5517.6
$ws.Range("J134").Value = 2601.3125
$ws.Range("K134").Value = 16552.8
$ws.Range("L134").Value = 7803.9375
$ws.Range("M134").Value = -14017.8
$ws.Range("N134").Value = -12873.9375
$ws.Range("H136").Value = 1470.1333
$ws.Range("I136").Value = 1426.9595
$ws.Range("K136").Value = 4280.8785
$ws.Range("M136").Value = -1730.8785

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4542.5
$ws.Range("I56").Value = 4542.5
$ws.Range("K56").Value = 4542.5
$ws.Range("M56").Value = -4012.5
$ws.Range("H107").Value = 62973.688
$ws.Range("I107").Value = 465.27274
$ws.Range("J107").Value = 200492.2
$ws.Range("K107").Value = 1395.81822
$ws.Range("L107").Value = 601476.6000000001
$ws.Range("M107").Value = 524.1817799999999
$ws.Range("N107").Value = -605316.6000000001
$ws.Range("H113").Value = 778.8570999999999
$ws.Range("I113").Value = 790
$ws.Range("K113").Value = 2370
$ws.Range("M113").Value = -200
$ws.Range("H131").Value = 6757580
$ws.Range("J131").Value = 856.4143
$ws.Range("L131").Value = 2569.2429
$ws.Range("N131").Value = -12649.2429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5931.899
$ws.Range("I70").Value = 5610.547
$ws.Range("K70").Value = 5610.547
$ws.Range("M70").Value = -5340.547
$ws.Range("H73").Value = 5931.899
$ws.Range("I73").Value = 5610.547
$ws.Range("K73").Value = 5610.547
$ws.Range("M73").Value = -4674.547
$ws.Range("H136").Value = 16326.214
$ws.Range("J136").Value = 16326.214
$ws.Range("L136").Value = 48978.642
$ws.Range("N136").Value = -54078.642
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2838.6667
$ws.Range("I7").Value = 1619.5
$ws.Range("J7").Value = 6496.1665
$ws.Range("K7").Value = 1619.5
$ws.Range("L7").Value = 6496.1665
$ws.Range("M7").Value = -1507.5
$ws.Range("N7").Value = -6720.1665
$ws.Range("H16").Value = 653.75
$ws.Range("I16").Value = 653.75
$ws.Range("K16").Value = 653.75
$ws.Range("M16").Value = -483.75
$ws.Range("H122").Value = 3250.6428
$ws.Range("I122").Value = 1763
$ws.Range("K122").Value = 5289
$ws.Range("M122").Value = -2839
$ws.Range("H126").Value = 2838.6667
$ws.Range("I126").Value = 1619.5
$ws.Range("J126").Value = 6496.1665
$ws.Range("K126").Value = 4858.5
$ws.Range("L126").Value = 19488.4995
$ws.Range("M126").Value = -2388.5
$ws.Range("N126").Value = -24428.4995
$ws.Range("H132").Value = 5409.5654
$ws.Range("I132").Value = 1920.6364
$ws.Range("J132").Value = 8607.75
$ws.Range("K132").Value = 5761.9092
$ws.Range("L132").Value = 25823.25
$ws.Range("M132").Value = -3231.9092
$ws.Range("N132").Value = -30883.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 247
$ws.Range("I100").Value = 234.125
$ws.Range("K100").Value = 468.25
$ws.Range("M100").Value = 72.75
$ws.Range("H136").Value = 5851
$ws.Range("I136").Value = 4302.75
$ws.Range("K136").Value = 12908.25
$ws.Range("M136").Value = -10358.25
